{"js": "// Replace the \"Unraveling the Enigma of Consciousness\" essay with the\n// \"Symphony of Growth and Adaptation\" biology essay, including the new\n// author / email, new body text, and a trailing empty paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Expected original layout (7 paragraphs):\n//   0: Title\n//   1: Author name\n//   2: Email address\n//   3: (empty spacer paragraph)\n//   4: Main essay body (contains manual line breaks \\v)\n//   5: \"Summary\" heading\n//   6: Summary paragraph\nconst titlePara = items[0];\nconst authorPara = items[1];\nconst emailPara = items[2];\nconst bodyPara = items[4];\nconst summaryParaBody = items[6];\n\nconst newTitle =\n  \"The Symphony of Growth and Adaptation: Exploring the Wonders of Biology\";\nconst newAuthor = \"Myra Whitaker\";\nconst newEmail = \"myra.whitaker399@schoolmail.edu\";\n\nconst newBody =\n  \"Biology, the study of life, embarks us on an enthralling journey into the intricate world of living organisms.\" +\n  \" Imagine yourself as an orchestra conductor, delicately orchestrating the symphony of life's processes.\" +\n  \" Cells, the fundamental units of life, are like harmonious instruments, each contributing its unique melody to the overall composition.\" +\n  \" From the smallest microbes to the vast blue whale, each organism dances to its own rhythm, adapting to the ever-changing notes of the environment.\" +\n  \" As you delve into the pages of biology, you'll unveil the secrets of life's intricate score, a symphony of growth, adaptation, and resilience.\" +\n  \"\\v\\v\" +\n  \"The beauty of biology lies in its ability to connect us to the world around us.\" +\n  \" It's a science that weaves together the tapestry of life, unraveling the threads of complex systems and unveiling the interdependence of all organisms.\" +\n  \" From the delicate dance of pollination to the intricate food webs that sustain ecosystems, biology showcases the interconnectedness of all living things.\" +\n  \" As you embark on this journey, you'll explore the marvels of biodiversity, discovering the diverse adaptations that allow creatures to thrive in various environments.\" +\n  \"\\v\\v\" +\n  \"Biology, however, is not merely a collection of facts and theories; it's an exploration of the profound questions that have captivated humanity for eons.\" +\n  \" How did life begin? How do organisms evolve? How can we use our knowledge of biology to address global challenges such as disease, hunger, and climate change? As you navigate the depths of biology, you'll grapple with these questions, engaging in scientific inquiry and critical thinking to unravel the mysteries that lie ahead.\";\n\nconst newSummary =\n  \"Biology, the study of life, invites us to explore the captivating symphony of growth and adaptation.\" +\n  \" Through the study of cells, organisms, and ecosystems, we gain an understanding of life's intricate processes and our interconnectedness with the natural world.\" +\n  \" Biology challenges us to ask profound questions about the origins and evolution of life and empowers us to seek solutions to global challenges.\" +\n  \" As you delve into the realm of biology, you will not only acquire knowledge but also cultivate a sense of wonder and appreciation for the symphony of life that surrounds us.\";\n\ntitlePara.insertText(newTitle, \"Replace\");\nauthorPara.insertText(newAuthor, \"Replace\");\nemailPara.insertText(newEmail, \"Replace\");\nbodyPara.insertText(newBody, \"Replace\");\nsummaryParaBody.insertText(newSummary, \"Replace\");\n\n// A new empty paragraph is appended after the summary paragraph (before\n// the section break / end of the document body).\nsummaryParaBody.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# Replace the \"Unraveling the Enigma of Consciousness\" essay with the\n# \"Symphony of Growth and Adaptation\" biology essay: new title, new\n# author / email, new body paragraphs, and a trailing empty paragraph.\n\n$d = $word.ActiveDocument\n\n# Word's COM Range object returned directly from Paragraphs.Item(n).Range\n# only overwrites the first run when a paragraph has several runs, so we\n# rebuild a fresh Range over the same [Start,End) span before assigning\n# .Text -- that reliably replaces the whole paragraph (collapsing it to a\n# single run) while preserving the trailing paragraph mark.\nfunction Set-ParagraphText($paraIndex, $text) {\n    $p = $d.Paragraphs.Item($paraIndex)\n    $rng = $d.Range($p.Range.Start, $p.Range.End)\n    $rng.Text = $text\n}\n\n$vt = [char]11\n\n$newTitle = \"The Symphony of Growth and Adaptation: Exploring the Wonders of Biology\"\n$newAuthor = \"Myra Whitaker\"\n$newEmail = \"myra.whitaker399@schoolmail.edu\"\n\n$newBody = \"Biology, the study of life, embarks us on an enthralling journey into the intricate world of living organisms.\" `\n  + \" Imagine yourself as an orchestra conductor, delicately orchestrating the symphony of life's processes.\" `\n  + \" Cells, the fundamental units of life, are like harmonious instruments, each contributing its unique melody to the overall composition.\" `\n  + \" From the smallest microbes to the vast blue whale, each organism dances to its own rhythm, adapting to the ever-changing notes of the environment.\" `\n  + \" As you delve into the pages of biology, you'll unveil the secrets of life's intricate score, a symphony of growth, adaptation, and resilience.\" `\n  + $vt + $vt `\n  + \"The beauty of biology lies in its ability to connect us to the world around us.\" `\n  + \" It's a science that weaves together the tapestry of life, unraveling the threads of complex systems and unveiling the interdependence of all organisms.\" `\n  + \" From the delicate dance of pollination to the intricate food webs that sustain ecosystems, biology showcases the interconnectedness of all living things.\" `\n  + \" As you embark on this journey, you'll explore the marvels of biodiversity, discovering the diverse adaptations that allow creatures to thrive in various environments.\" `\n  + $vt + $vt `\n  + \"Biology, however, is not merely a collection of facts and theories; it's an exploration of the profound questions that have captivated humanity for eons.\" `\n  + \" How did life begin? How do organisms evolve? How can we use our knowledge of biology to address global challenges such as disease, hunger, and climate change? As you navigate the depths of biology, you'll grapple with these questions, engaging in scientific inquiry and critical thinking to unravel the mysteries that lie ahead.\"\n\n$newSummary = \"Biology, the study of life, invites us to explore the captivating symphony of growth and adaptation.\" `\n  + \" Through the study of cells, organisms, and ecosystems, we gain an understanding of life's intricate processes and our interconnectedness with the natural world.\" `\n  + \" Biology challenges us to ask profound questions about the origins and evolution of life and empowers us to seek solutions to global challenges.\" `\n  + \" As you delve into the realm of biology, you will not only acquire knowledge but also cultivate a sense of wonder and appreciation for the symphony of life that surrounds us.\"\n\n# Paragraph layout in the original document:\n#   1: Title\n#   2: Author name\n#   3: Email address\n#   4: (empty spacer paragraph)\n#   5: Main essay body (contains manual line breaks)\n#   6: \"Summary\" heading\n#   7: Summary paragraph\nSet-ParagraphText 1 $newTitle\nSet-ParagraphText 2 $newAuthor\nSet-ParagraphText 3 $newEmail\nSet-ParagraphText 5 $newBody\nSet-ParagraphText 7 $newSummary\n\n# A new empty paragraph is appended after the summary paragraph (before\n# the section break / end of the document body).\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n"}
